# 727-MS-EPP-DB-SAR-REC-NON-RNI-CTRFD-SAR-MD-TR-1-EarlyRePayment-Loanproduct.xlsx
# "Added periodic & upfront related scenarios"
#
# The "repaymentstrategy" input on the ProductLoanInput sheet (cell B17)
# is changed from "Mifos style" to "Penalties, Fees, Interest, Principal
# order", and the cell is (re)formatted with left/top alignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$cell = $ws.Range("B17")
$cell.Value = "Penalties, Fees, Interest, Principal order"

# Apply left/top alignment to the updated cell (xlLeft = -4131, xlTop = -4160)
$cell.HorizontalAlignment = -4131
$cell.VerticalAlignment = -4160

# Move/leave the selection on the cell that was just edited.
$cell.Select()
